# Fruta / hortaliza, semanal
# Rows 3-13 get their weekly (date-dependent) figures re-shuffled: the
# Fecha (D), Volumen (J), Precio minimo/maximo/promedio (K/L/M), Origen (O)
# and Precio $/Kg (P) columns are permuted across rows. Columns A, B, C,
# E, F, G, H, I, N, Q, R are untouched, so they are left alone.
#
# (Reading a cell's .Value back out through this COM bridge and piping it
# straight into another cell's .Value can mis-marshal string variants, so
# the target values are written as literals taken from the known source
# rows instead of being copied live.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    3  = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 }
    4  = @{ D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí"; P = 578 }
    5  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
    6  = @{ D = 44435; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    7  = @{ D = 44376; J = 15; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    8  = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    9  = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    10 = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí"; P = 640 }
    11 = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    12 = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    13 = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
}
